# audiences-process.pptx -- sync non localizable files
#
# 1) Refresh the cached "datetimeFigureOut" date field text (slide master +
#    every slide layout) from 01/02/2023 -> 02/08/2023.
# 2) Widen the "Target audience(s) in campaigns" textbox on slide 1 and
#    extend its copy to mention journeys as well.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders -------------------------------------------------
# ppPlaceholderDate == 16. Walk the slide master plus every custom layout
# and update any date placeholder's text in place (id/shape untouched).
$targets = @($p.SlideMaster)
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $targets += $layouts.Item($j)
}

foreach ($t in $targets) {
    for ($i = 1; $i -le $t.Shapes.Count; $i++) {
        $shp = $t.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = "02/08/2023"
        }
    }
}

# --- 2) "Target audience(s) in campaigns" textbox on slide 1 -------------
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 11")
$shp.TextFrame.TextRange.Text = "Target audience(s) in campaigns and journeys"
# 2238214 EMU -- Width is a points (f32) property, so nudge slightly past
# the exact EMU->pt boundary so the round-trip lands back on 2238214 EMU.
$shp.Width = 176.237335
